# ESG_Dashboard.xlsx edit script
# - Fix "Production sync" baseline figures on STRATEGY SELECTOR.
# - Fix broken cross-sheet references that used underscore pseudo-names
#   (IMPACT_CONFIG / STRATEGY_SELECTOR) instead of the real, space-containing
#   sheet names ("IMPACT CONFIG" / "STRATEGY SELECTOR"), which must be
#   single-quoted in A1 formulas.
# - Overtime logic: new rows on CROSS REFERENCE, plus updated Logistics figures.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# helper: write a value as literal TEXT without disturbing the cell's
# existing style/number-format (typing a numeric-looking string into
# Range.Value auto-coerces it to a Number and can drag in the precedent's
# number format when used with .Formula; round-tripping a text formula
# through PasteSpecial(values) keeps both the original style AND text type).
# ---------------------------------------------------------------------------
function Set-TextValue($range, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)   # xlPasteValues
}

# ---------------------------------------------------------------------------
# 1) STRATEGY SELECTOR - baseline inputs + broken sheet references
# ---------------------------------------------------------------------------
$strategy = $wb.Worksheets.Item("STRATEGY SELECTOR")

$strategy.Range("B6").Value = 15000
$strategy.Range("B8").Value = 1200000

$strategy.Range("B7").Formula  = "=B6*'IMPACT CONFIG'!`$B`$4"

$strategy.Range("C13").Formula = "=B13*'IMPACT CONFIG'!B8"
$strategy.Range("D13").Formula = "=B13*'IMPACT CONFIG'!C8"
$strategy.Range("E13").Formula = "=D13*'IMPACT CONFIG'!`$B`$4"
$strategy.Range("I13").Formula = "='IMPACT CONFIG'!`$B`$4"
$strategy.Range("I13").ClearFormats()

$strategy.Range("C14").Formula = "=B14*'IMPACT CONFIG'!B9"
$strategy.Range("D14").Formula = "=B14*'IMPACT CONFIG'!C9"
$strategy.Range("E14").Formula = "=D14*'IMPACT CONFIG'!`$B`$4"
$strategy.Range("I14").Formula = "='IMPACT CONFIG'!`$B`$4"
$strategy.Range("I14").ClearFormats()

$strategy.Range("C15").Formula = "=`$B`$8*B15*'IMPACT CONFIG'!B10"
$strategy.Range("D15").Formula = "=`$B`$8*B15*'IMPACT CONFIG'!C10"
$strategy.Range("E15").Formula = "=D15*'IMPACT CONFIG'!`$B`$4"
$strategy.Range("I15").Formula = "='IMPACT CONFIG'!`$B`$4"
$strategy.Range("I15").ClearFormats()

$strategy.Range("C16").Formula = "=B16*'IMPACT CONFIG'!B11"
$strategy.Range("D16").Formula = "=B16*'IMPACT CONFIG'!C11"
$strategy.Range("E16").Formula = "=D16*'IMPACT CONFIG'!`$B`$4"
$strategy.Range("I16").Formula = "='IMPACT CONFIG'!`$B`$4"
$strategy.Range("I16").ClearFormats()

$strategy.Range("B22").Formula = "=MAX(0,(B6-B20)*'IMPACT CONFIG'!B4)"

# ---------------------------------------------------------------------------
# 2) UPLOAD READY ESG - broken sheet references (all target cells start
#    out unstyled, so ClearFormats() afterwards restores that "no explicit
#    style" state in case a precedent's number format got auto-applied).
# ---------------------------------------------------------------------------
$uploadReady = $wb.Worksheets.Item("UPLOAD READY ESG")

$uploadReady.Range("B6").Formula = "='STRATEGY SELECTOR'!B13"
$uploadReady.Range("C6").Formula = "='STRATEGY SELECTOR'!C13"
$uploadReady.Range("D6").Formula = "='STRATEGY SELECTOR'!D13"
$uploadReady.Range("E6").Formula = "='STRATEGY SELECTOR'!E13"

$uploadReady.Range("B7").Formula = "='STRATEGY SELECTOR'!B14"
$uploadReady.Range("C7").Formula = "='STRATEGY SELECTOR'!C14"
$uploadReady.Range("D7").Formula = "='STRATEGY SELECTOR'!D14"
$uploadReady.Range("E7").Formula = "='STRATEGY SELECTOR'!E14"

$uploadReady.Range("B8").Formula = "='STRATEGY SELECTOR'!B15"
$uploadReady.Range("C8").Formula = "='STRATEGY SELECTOR'!C15"
$uploadReady.Range("D8").Formula = "='STRATEGY SELECTOR'!D15"
$uploadReady.Range("E8").Formula = "='STRATEGY SELECTOR'!E15"

$uploadReady.Range("B9").Formula = "='STRATEGY SELECTOR'!B16"
$uploadReady.Range("C9").Formula = "='STRATEGY SELECTOR'!C16"
$uploadReady.Range("D9").Formula = "='STRATEGY SELECTOR'!D16"
$uploadReady.Range("E9").Formula = "='STRATEGY SELECTOR'!E16"

$uploadReady.Range("B6:E9").ClearFormats()

# ---------------------------------------------------------------------------
# 3) CROSS REFERENCE - Overtime/Unit Cost/Sales rows + Logistics update
# ---------------------------------------------------------------------------
$crossRef = $wb.Worksheets.Item("CROSS REFERENCE")

# Total Production figure resets to 0 pending new production sync
Set-TextValue $crossRef.Range("B5") "0"

# Insert three new rows ahead of the old "Logistics (Transport)" block
# (old rows 9-11 become rows 12-14), and carry the existing row-6 cell
# formatting (border/fill/font) down onto them so the new rows visually
# match the rest of the table.
$crossRef.Range("A7:B9").EntireRow.Insert()
$crossRef.Range("A6:B6").Copy()
$crossRef.Range("A7:B9").PasteSpecial(-4122)

$crossRef.Range("A7").Value = "Overtime Hours"
Set-TextValue $crossRef.Range("B7") "0"

$crossRef.Range("A8").Value = "Unit Cost Avg"
Set-TextValue $crossRef.Range("B8") "`$40.00"

$crossRef.Range("A9").Value = "Total Sales (Target)"
Set-TextValue $crossRef.Range("B9") "0"

# Updated Logistics Costs figure (now at row 13 after the insert above)
Set-TextValue $crossRef.Range("B13") "`$125,000"
